$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 142 (pushes existing rows 142..257 down to 143..258)
$ws.Rows("142:142").Insert()

# Populate the new row 142 with the new "Especial" / Provincia de Melipilla record
$ws.Range("A142").Value = 10
$ws.Range("B142").Value = "Vega Modelo de Temuco"
$ws.Range("C142").Value = "La Araucanía"
$ws.Range("D142").Value = 44827
$ws.Range("E142").Value = 9
$ws.Range("F142").Value = "Fruta"
$ws.Range("G142").Value = 100101
$ws.Range("H142").Value = "Berries"
$ws.Range("I142").Value = 100112025
$ws.Range("J142").Value = "Frutilla"
$ws.Range("K142").Value = "Sin especificar"
$ws.Range("L142").Value = "Especial"
$ws.Range("M142").Value = 45
$ws.Range("N142").Value = 18000
$ws.Range("O142").Value = 18000
$ws.Range("P142").Value = 18000
$ws.Range("Q142").Value = "$/bandeja 7 kilos"
$ws.Range("R142").Value = "Provincia de Melipilla"
$ws.Range("S142").Value = 2571
$ws.Range("T142").Value = 7
